# GuiaDidactica_LE_05_05.docx — "Versión final post Corrección de estilo"
#
# All substantive edits are pure copy-editing tweaks. Per the target
# revision, the *removed* characters are folded into the base text
# directly (no tracked deletions appear anywhere in the final XML),
# while the *newly added* characters are recorded as tracked insertions
# (w:ins) authored by "Admincmovil". There is also a relocation of the
# lone "_GoBack" bookmark, and a cosmetic run-merge that drops stale
# w:proofErr spell-check markers around "presaberes".

$d = $word.ActiveDocument
$word.UserName = "Admincmovil"

function Get-Text {
    return $d.Content.Text
}

# ---------------------------------------------------------------
# Edit 1: "Para lo cual," -> "Para lo cual" + ins(":")
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("Para lo cual,")
$delStart = $idx + "Para lo cual".Length
$delRng = $d.Range($delStart, $delStart + 1)
$delRng.Delete()

$d.TrackRevisions = $true
$insPos = $d.Range($delStart, $delStart)
$insPos.InsertAfter(":")

# ---------------------------------------------------------------
# Edit 2: "...verbos regulares e los irregulares." -> "...regulares " + ins("y ") + "los irregulares."
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("regulares e los irregulares")
$delStart = $idx + "regulares ".Length
$delRng = $d.Range($delStart, $delStart + 2)
$delRng.Delete()

$d.TrackRevisions = $true
$insPos = $d.Range($delStart, $delStart)
$insPos.InsertAfter("y ")

# ---------------------------------------------------------------
# Edit 3: "capacidad humana por excelencia" -> "capacidad humana" + ins(",") + " por excelencia" + ins(",")
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("capacidad humana por excelencia")
$pos1 = $idx + "capacidad humana".Length

$d.TrackRevisions = $true
$insPos1 = $d.Range($pos1, $pos1)
$insPos1.InsertAfter(",")

$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("capacidad humana, por excelencia")
$pos2 = $idx + "capacidad humana, por excelencia".Length

$d.TrackRevisions = $true
$insPos2 = $d.Range($pos2, $pos2)
$insPos2.InsertAfter(",")

# ---------------------------------------------------------------
# Edit 4: "...y concretar sus procesos..." -> "...y concretar" + ins("á") + " sus procesos..."
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("concretar sus procesos de significación")
$pos = $idx + "concretar".Length

$d.TrackRevisions = $true
$insPos = $d.Range($pos, $pos)
$insPos.InsertAfter("á")

# ---------------------------------------------------------------
# Edit 5: "...dispuestas aquí para la estudio del verbo..." -> "...para " + ins("el ") + "estudio del verbo..."
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("aquí para la estudio del verbo")
$delStart = $idx + "aquí para ".Length
$delRng = $d.Range($delStart, $delStart + 3)
$delRng.Delete()

$d.TrackRevisions = $true
$insPos = $d.Range($delStart, $delStart)
$insPos.InsertAfter("el ")

# ---------------------------------------------------------------
# Edit 6: "...a su disposición tanto los instrumentos..." -> "...disposición" + ins(",") + " tanto los instrumentos..."
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("a su disposición tanto los instrumentos")
$pos = $idx + "a su disposición".Length

$d.TrackRevisions = $true
$insPos = $d.Range($pos, $pos)
$insPos.InsertAfter(",")

# ---------------------------------------------------------------
# Edit 7: "...recursos ofrecidos y la habilidad..." -> "...ofrecidos" + ins(",") + " y la habilidad..."
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("recursos ofrecidos y la habilidad")
$pos = $idx + "recursos ofrecidos".Length

$d.TrackRevisions = $true
$insPos = $d.Range($pos, $pos)
$insPos.InsertAfter(",")

# ---------------------------------------------------------------
# Edit 8: relocate the lone "_GoBack" bookmark from after "práctica"
# to between "la mediación entre " and "este y los presaberes".
# Adding a bookmark named "_GoBack" elsewhere moves the existing one.
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("la mediación entre este y los")
$pos = $idx + "la mediación entre ".Length
$bmRng = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------
# Edit 9: drop the stale proofErr spell-check wrapper around "presaberes"
# by touching (delete + reinsert) the whole sentence; the engine then
# re-merges the three runs into one, matching the cleaned-up markup.
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$needle = "sus presaberes y nuevos conocimientos al momento de conjugar verbos en diferentes modos, tiempos y personas."
$idx = $full.IndexOf($needle)
$rng = $d.Range($idx, $idx + $needle.Length)
$rng.Delete()
$insPos = $d.Range($idx, $idx)
$insPos.InsertAfter($needle)

# ---------------------------------------------------------------
# Edit 10: "...tanto de análisis como de síntesis..." ->
#   "...tanto de" + ins("l") + " análisis como de" + ins(" la") + " síntesis..."
# ---------------------------------------------------------------
$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("tanto de análisis como de síntesis")
$pos1 = $idx + "tanto de".Length

$d.TrackRevisions = $true
$insPos1 = $d.Range($pos1, $pos1)
$insPos1.InsertAfter("l")

$d.TrackRevisions = $false
$full = Get-Text
$idx = $full.IndexOf("tanto del análisis como de síntesis")
$pos2 = $idx + "tanto del análisis como de".Length

$d.TrackRevisions = $true
$insPos2 = $d.Range($pos2, $pos2)
$insPos2.InsertAfter(" la")

$d.TrackRevisions = $false

Write-Output "edits applied"
